$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New SKU rows to append (SKU, product weight, package weight).
$newRows = @(
    @("MBK XXL", 115, 13),
    @("MBK XL",  115, 13),
    @("MBK L",   127, 13),
    @("MBK M",   137, 13),
    @("MBK S",   145, 13),
    @("MBK XS",  145, 13),
    @("MBV XL",  160, 13),
    @("MBV L",   115, 13),
    @("MBV M",   127, 13),
    @("MBV S",   137, 13),
    @("HO XL",   145, 13),
    @("HO L",    160, 13),
    @("HO M",    125, 13),
    @("HO S",    129, 13),
    @("SW XL",   115, 13),
    @("SW L",    127, 13),
    @("SW M",    137, 13),
    @("SW S",    145, 13),
    @("MPM XL",  160, 13),
    @("MPM L",   115, 13),
    @("MPM M",   127, 13),
    @("MPM S",   137, 13),
    @("MBM XL",  145, 13),
    @("MBM L",   160, 13),
    @("MBM M",   115, 13),
    @("MBM S",   127, 13),
    @("MPV XXL", 160, 13),
    @("MPV XL",  160, 13),
    @("MPV L",   115, 13),
    @("MPV M",   160, 13),
    @("MPV S",   115, 13),
    @("MPV XS",  115, 13)
)

$startRow = 69
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

$ws.Range("E73").Select()
